# Insert a new weekly price record for "Feria Lagunitas de Puerto Montt -
# Poroto verde" above the current row 114, shifting the existing rows
# 114-130 down to 115-131 (the sheet dimension grows from A1:R130 to
# A1:R131).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above row 114; Excel pushes rows 114:130 down to
# 115:131 and copies row-114's formatting into the new row (same as a
# manual "Insert Sheet Rows" above the selection).
$ws.Rows.Item(114).Insert()

# Populate the freshly inserted row 114 with the new weekly observation.
$ws.Range("A114").Value = 4
$ws.Range("B114").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C114").Value = "Los Lagos"
$ws.Range("D114").Value = 44995
$ws.Range("E114").Value = 10
$ws.Range("F114").Value = 100112031
$ws.Range("G114").Value = "Poroto verde"
$ws.Range("H114").Value = "Magnum"
$ws.Range("I114").Value = "Primera"
$ws.Range("J114").Value = 40
$ws.Range("K114").Value = 30000
$ws.Range("L114").Value = 30000
$ws.Range("M114").Value = 30000
$ws.Range("N114").Value = "$/saco 25 kilos"
$ws.Range("O114").Value = "Región Metropolitana"
$ws.Range("P114").Value = 1200
$ws.Range("Q114").Value = 25
$ws.Range("R114").Value = "Hortaliza"
